$p = $ppt.ActivePresentation
Write-Host $p.Slides.Count
